$d = $word.ActiveDocument

# 1. "we made assigned these" -> "we assigned these"
$d.Content.Find.Execute(
    "we made assigned these to our group members",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "we assigned these to our group members", 2)

# 2. "LED distance sensor: " -> "LED antenna sensor: "
$d.Content.Find.Execute(
    "LED distance sensor: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LED antenna sensor: ", 2)

# 3. "Maximum speed of the mock up robot is 12 cm/s." -> add trailing note
$d.Content.Find.Execute(
    "Maximum speed of the mock up robot is 12 cm/s.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Maximum speed of the mock up robot is 12 cm/s. (previously it was 20cm/s)", 2)

# 4. "Distance between the center of the robot and ground..." -> insert clause
$d.Content.Find.Execute(
    "Distance between the center of the robot and ground needs to in between 12 cm and 16 cm.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Distance between the center of the rectangle (visibility marker) on the robot and ground needs to in between 12 cm and 16 cm.", 2)
